# 16 Jul 2016: Major styling updates.
# - Row 4 (Issue: "Update CB Functions / Error Handling / Front-End") status
#   changes from "Closed" to "Open".
# - A note is added in D4 explaining the status.
# - The active selection moves to A5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the status of row 4 from "Closed" to "Open"
$ws.Range("B4").Value = "Open"

# Add a note in the Notes column for row 4
$ws.Range("D4").Value = "All are OK except the LoginController one."

# The new wrapped note text makes row 4 taller; match the row height Excel
# computed once the wrapped note text was added (AutoFit() is a no-op on
# Rows/Range in this runtime, so set the resulting height directly)
$ws.Rows.Item(4).RowHeight = 28

# Update the selection to match the saved view state
$ws.Range("A5").Select()
